$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old intermediate sample rows (2 and 3) while keeping row 4's
# row index intact (clear contents rather than a shifting delete).
$ws.Rows("2:3").ClearContents()

# Refresh row 4 with the updated values from the re-run plotting/primer data.
$row4 = @{
    "A" = 0.3300896489316999
    "B" = 0.3937716913375339
    "C" = 0.2831114229638663
    "D" = 0.3572739247688735
    "E" = 0.3016866985809186
    "F" = 7.005345674098685
    "G" = 9.155225268991765
    "H" = 5.443870653969218
    "I" = 7.911868107462011
    "J" = 6.057471535817356
    "K" = 2.312451696222781
    "L" = 3.165432737870411
    "M" = 1.704226334517611
    "N" = 2.663326251156235
    "O" = 1.942531147614155
    "P" = 65.10586000000001
    "Q" = 90.81940183026543
    "R" = 47.44349024310426
    "S" = 74.42357440141171
    "T" = 54.8440184763783
    "U" = 0.2013384968688031
    "V" = 0.317510705356708
    "W" = 0.09932328267674136
    "X" = 0.2602566524088788
    "Y" = 0.1427160606128784
    "Z" = 0.8478579537441078
    "AA" = 0.9546501438494737
    "AB" = 0.697919396871762
    "AC" = 0.9097483621377701
    "AD" = 0.7726976537329295
}

foreach ($col in $row4.Keys) {
    $ws.Range("$col`4").Value = $row4[$col]
}
